# Update the addition/subtraction practice sheet: replace each of the
# 100 equation cells (5 columns x 20 rows) in the single table with its
# new value, matching cell-by-cell so duplicate old values (e.g. two
# cells both reading "30+30=60") resolve to the correct distinct
# replacement instead of a global text search/replace.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "67-20=47"
$t.Cell(1, 2).Range.Text = "61-60=1"
$t.Cell(1, 3).Range.Text = "67+22=89"
$t.Cell(1, 4).Range.Text = "63-38=25"
$t.Cell(1, 5).Range.Text = "10+50=60"

$t.Cell(2, 1).Range.Text = "27+16=43"
$t.Cell(2, 2).Range.Text = "52-20=32"
$t.Cell(2, 3).Range.Text = "87-20=67"
$t.Cell(2, 4).Range.Text = "64+34=98"
$t.Cell(2, 5).Range.Text = "2+41=43"

$t.Cell(3, 1).Range.Text = "64-1=63"
$t.Cell(3, 2).Range.Text = "65+0=65"
$t.Cell(3, 3).Range.Text = "1+45=46"
$t.Cell(3, 4).Range.Text = "20+40=60"
$t.Cell(3, 5).Range.Text = "57+12=69"

$t.Cell(4, 1).Range.Text = "25+13=38"
$t.Cell(4, 2).Range.Text = "43+29=72"
$t.Cell(4, 3).Range.Text = "86-4=82"
$t.Cell(4, 4).Range.Text = "83-51=32"
$t.Cell(4, 5).Range.Text = "96-50=46"

$t.Cell(5, 1).Range.Text = "21+35=56"
$t.Cell(5, 2).Range.Text = "25+31=56"
$t.Cell(5, 3).Range.Text = "8+43=51"
$t.Cell(5, 4).Range.Text = "8-0=8"
$t.Cell(5, 5).Range.Text = "24+22=46"

$t.Cell(6, 1).Range.Text = "94-72=22"
$t.Cell(6, 2).Range.Text = "31-30=1"
$t.Cell(6, 3).Range.Text = "46+38=84"
$t.Cell(6, 4).Range.Text = "15+23=38"
$t.Cell(6, 5).Range.Text = "99-74=25"

$t.Cell(7, 1).Range.Text = "31-22=9"
$t.Cell(7, 2).Range.Text = "84-36=48"
$t.Cell(7, 3).Range.Text = "62-43=19"
$t.Cell(7, 4).Range.Text = "57-44=13"
$t.Cell(7, 5).Range.Text = "87-29=58"

$t.Cell(8, 1).Range.Text = "55-8=47"
$t.Cell(8, 2).Range.Text = "53-21=32"
$t.Cell(8, 3).Range.Text = "13+46=59"
$t.Cell(8, 4).Range.Text = "48-38=10"
$t.Cell(8, 5).Range.Text = "8+12=20"

$t.Cell(9, 1).Range.Text = "16-9=7"
$t.Cell(9, 2).Range.Text = "92-83=9"
$t.Cell(9, 3).Range.Text = "20+5=25"
$t.Cell(9, 4).Range.Text = "41+26=67"
$t.Cell(9, 5).Range.Text = "40-15=25"

$t.Cell(10, 1).Range.Text = "96-87=9"
$t.Cell(10, 2).Range.Text = "42-36=6"
$t.Cell(10, 3).Range.Text = "16+78=94"
$t.Cell(10, 4).Range.Text = "42-16=26"
$t.Cell(10, 5).Range.Text = "10+3=13"

$t.Cell(11, 1).Range.Text = "39+44=83"
$t.Cell(11, 2).Range.Text = "86-15=71"
$t.Cell(11, 3).Range.Text = "44+41=85"
$t.Cell(11, 4).Range.Text = "63+9=72"
$t.Cell(11, 5).Range.Text = "29+70=99"

$t.Cell(12, 1).Range.Text = "8+15=23"
$t.Cell(12, 2).Range.Text = "26+37=63"
$t.Cell(12, 3).Range.Text = "95-58=37"
$t.Cell(12, 4).Range.Text = "45-10=35"
$t.Cell(12, 5).Range.Text = "45+3=48"

$t.Cell(13, 1).Range.Text = "89-47=42"
$t.Cell(13, 2).Range.Text = "68+5=73"
$t.Cell(13, 3).Range.Text = "3+8=11"
$t.Cell(13, 4).Range.Text = "55-44=11"
$t.Cell(13, 5).Range.Text = "14+23=37"

$t.Cell(14, 1).Range.Text = "38+33=71"
$t.Cell(14, 2).Range.Text = "55-0=55"
$t.Cell(14, 3).Range.Text = "82+9=91"
$t.Cell(14, 4).Range.Text = "58+37=95"
$t.Cell(14, 5).Range.Text = "43+6=49"

$t.Cell(15, 1).Range.Text = "23+30=53"
$t.Cell(15, 2).Range.Text = "55-36=19"
$t.Cell(15, 3).Range.Text = "36-34=2"
$t.Cell(15, 4).Range.Text = "18+48=66"
$t.Cell(15, 5).Range.Text = "44+25=69"

$t.Cell(16, 1).Range.Text = "85-1=84"
$t.Cell(16, 2).Range.Text = "86+3=89"
$t.Cell(16, 3).Range.Text = "27+6=33"
$t.Cell(16, 4).Range.Text = "95-87=8"
$t.Cell(16, 5).Range.Text = "58+14=72"

$t.Cell(17, 1).Range.Text = "69-65=4"
$t.Cell(17, 2).Range.Text = "80+9=89"
$t.Cell(17, 3).Range.Text = "56-25=31"
$t.Cell(17, 4).Range.Text = "34+37=71"
$t.Cell(17, 5).Range.Text = "65+2=67"

$t.Cell(18, 1).Range.Text = "8+4=12"
$t.Cell(18, 2).Range.Text = "37-14=23"
$t.Cell(18, 3).Range.Text = "45+38=83"
$t.Cell(18, 4).Range.Text = "12+26=38"
$t.Cell(18, 5).Range.Text = "16-14=2"

$t.Cell(19, 1).Range.Text = "85-53=32"
$t.Cell(19, 2).Range.Text = "39-10=29"
$t.Cell(19, 3).Range.Text = "16+82=98"
$t.Cell(19, 4).Range.Text = "46-39=7"
$t.Cell(19, 5).Range.Text = "88-9=79"

$t.Cell(20, 1).Range.Text = "67+14=81"
$t.Cell(20, 2).Range.Text = "66-45=21"
$t.Cell(20, 3).Range.Text = "26-22=4"
$t.Cell(20, 4).Range.Text = "48+45=93"
$t.Cell(20, 5).Range.Text = "41+17=58"

